$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last question (row 13) was saved without its "reference" link in
# column C. Copy the reference cell from the row above (C12) into C13 -
# this carries over both the value ("https://www.lanyingim.com/", reusing
# the existing shared string) and the cell's formatting/style.
$ws.Range("C12").Copy($ws.Range("C13"))

# Row 13 now holds a 14pt-font cell like the other data rows, so its
# height grows to match them.
$ws.Rows.Item(13).RowHeight = 19

# Move the active selection to the newly completed cell.
$ws.Range("C13").Select()
